# Apply "Use new string resources" changes to the Functions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Rows whose "Written" column (B) should now reference the "X" string resource.
# (Set cell values before hiding rows, so row heights are not recalculated.)
$rowsWithX = @(120, 121, 331, 373)
foreach ($r in $rowsWithX) {
    $ws.Cells.Item($r, 2).Value = "X"
}

# Rows that become hidden (they were previously visible).
$rowsToHide = @(89, 120, 121, 331, 373)
foreach ($r in $rowsToHide) {
    $ws.Rows.Item($r).Hidden = $true
}

# Update the active selection on the sheet to A269.
$ws.Range("A269").Select()
